$d = $word.ActiveDocument

function Find-Range($searchText) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r
}

function Insert-TaggedRun($pos, $tagText) {
    # Build a zero-length range at $pos and stamp it with the formatting of an
    # existing run that already contains the exact literal tag text (so the
    # new run picks up the identical rFonts/color/sz/szCs/rtl properties).
    $tmpl = Find-Range $tagText
    $ft = $tmpl.FormattedText
    $ins = $d.Range($pos, $pos)
    $ins.FormattedText = $ft
}

# ---------------------------------------------------------------------------
# Edit 1: "Pieds de petit" -> "Pieds de " + <al> + "petit"
# ---------------------------------------------------------------------------
$f = Find-Range "Pieds de petit"
$insertPos = $f.Start + "Pieds de ".Length
Insert-TaggedRun $insertPos "<al>"

# ---------------------------------------------------------------------------
# Edit 2: "es " + <al> + "lezardes" -> "es lezardes"
# ---------------------------------------------------------------------------
$f = Find-Range "es <al>lezardes"
$f.Text = "es lezardes"

# ---------------------------------------------------------------------------
# Edit 3: " puys gecte le second moule Et la " ->
#         " puys gecte le second " + <tl> + "moule" + </tl> + " Et la "
# ---------------------------------------------------------------------------
$f = Find-Range " puys gecte le second moule Et la "
$prefixLen = " puys gecte le second ".Length
$insertPos1 = $f.Start + $prefixLen
Insert-TaggedRun $insertPos1 "<tl>"

$f2 = Find-Range "moule Et la"
$insertPos2 = $f2.Start + "moule".Length
Insert-TaggedRun $insertPos2 "</tl>"

# ---------------------------------------------------------------------------
# Edit 4: "left-top" -> "left-" + "middle" (plain run, no color)
# ---------------------------------------------------------------------------
function Get-MiddleTemplate() {
    $t = Find-Range "left-middle"
    $midStart = $t.Start + "left-".Length
    $midRange = $d.Range($midStart, $t.End)
    return $midRange.FormattedText
}

$f = Find-Range "left-top"
$cutPos = $f.Start + "left-".Length
$tailRange = $d.Range($cutPos, $f.End)
$tailRange.Delete()
$ins = $d.Range($cutPos, $cutPos)
$ins.FormattedText = Get-MiddleTemplate

# ---------------------------------------------------------------------------
# Edit 5: "grands " + <al> + "lesards" -> <al> + "grands lesards"
# ---------------------------------------------------------------------------
$f = Find-Range "grands <al>lesards"
$grandsRange = $d.Range($f.Start, $f.Start + "grands ".Length)
$grandsRange.Delete()

$f2 = Find-Range "lesards"
$f2.Text = "grands lesards"

# ---------------------------------------------------------------------------
# Edit 6: "cire" + " dure ou " -> "cire dure" + " ou "
# ---------------------------------------------------------------------------
$f = Find-Range "<m>cire</m> dure ou "
$cireStart = $f.Start + "<m>".Length
$cireEnd = $cireStart + "cire".Length
$cireRange = $d.Range($cireStart, $cireEnd)
$cireRange.Text = "cire dure"

$f2 = Find-Range " dure ou "
$f2.Text = " ou "

# ---------------------------------------------------------------------------
# Edit 7: "</m></tl>" + " chault y faire tenir la gorge de lanimal" ->
#         "</m>" + " chault" + </tl> + " y faire tenir la gorge de lanimal"
# ---------------------------------------------------------------------------
$f = Find-Range "</m></tl> chault y faire tenir la gorge de lanimal"
$tagRange = $d.Range($f.Start, $f.Start + "</m></tl>".Length)
$tagRange.Text = "</m>"

$f2 = Find-Range " chault y faire tenir la gorge de lanimal"
$insertPos = $f2.Start + " chault".Length
Insert-TaggedRun $insertPos "</tl>"

# ---------------------------------------------------------------------------
# Edit 8: "left-bottom" -> "left-" + "middle" (plain run, no color)
# ---------------------------------------------------------------------------
$f = Find-Range "left-bottom"
$cutPos = $f.Start + "left-".Length
$tailRange = $d.Range($cutPos, $f.End)
$tailRange.Delete()
$ins = $d.Range($cutPos, $cutPos)
$ins.FormattedText = Get-MiddleTemplate

Write-Output "All edits applied."
